# Auto-generated Excel COM-interop edit script
# Applies the weekly CompStat data refresh (Volume 31 Number 19 -> 20,
# week of 5/13/2024-5/19/2024) plus the updated crime-stat figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number "19" -> "20" (edit in place, keep rich-text run) ---
$hdr = $ws.Range("A8")
$hdrVal = $hdr.Value()
$pos = $hdrVal.IndexOf("19") + 1
$hdr.Characters($pos, 2).Text = "20"

# --- Header: report week dates "5/6/2024"-"5/12/2024" -> "5/13/2024"-"5/19/2024" ---
$wk = $ws.Range("C9")
$wkVal = $wk.Value()
$posEnd = $wkVal.IndexOf("5/12/2024") + 1
$wk.Characters($posEnd, 9).Text = "5/19/2024"
$posStart = $wkVal.IndexOf("5/6/2024") + 1
$wk.Characters($posStart, 8).Text = "5/13/2024"

# --- Column H (8) width bump to match new best-fit ---
$ws.Columns.Item(8).ColumnWidth = 7.433768

# --- Data cells (rows 15-31) ---
$ws.Range("L15").Value = -25
$ws.Range("C16").Value = "'0"
$ws.Range("A16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D16").Value = "'0"
$ws.Range("A16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E16").Value = "'***.*"
$ws.Range("A16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 54
$ws.Range("K16").Value = 58.823529411764
$ws.Range("L16").Value = 25.581395348837
$ws.Range("M16").Value = -40
$ws.Range("N16").Value = -77.959183673469
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 133.333333333333
$ws.Range("I17").Value = 54
$ws.Range("J17").Value = 58
$ws.Range("K17").Value = -6.896551724137
$ws.Range("L17").Value = -14.285714285714
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = -51.351351351351
$ws.Range("C18").Value = 3
$ws.Range("F18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -11.111111111111
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = -39.024390243902
$ws.Range("L18").Value = -40.47619047619
$ws.Range("M18").Value = -74.226804123711
$ws.Range("N18").Value = -94.845360824742
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 7.142857142857
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -7.407407407407
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 234
$ws.Range("K19").Value = -1.709401709401
$ws.Range("L19").Value = 32.183908045977
$ws.Range("M19").Value = 26.373626373626
$ws.Range("N19").Value = -14.49814126394
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 183.333333333333
$ws.Range("I20").Value = 59
$ws.Range("J20").Value = 39
$ws.Range("K20").Value = 51.282051282051
$ws.Range("L20").Value = 43.90243902439
$ws.Range("M20").Value = -11.940298507462
$ws.Range("N20").Value = -94.722719141323
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 13.043478260869
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = 19.512195121951
$ws.Range("I21").Value = 425
$ws.Range("J21").Value = 412
$ws.Range("K21").Value = 3.155339805825
$ws.Range("L21").Value = 15.803814713896
$ws.Range("M21").Value = -12.551440329218
$ws.Range("N21").Value = -81.077471059661
$ws.Range("C23").Value = 2
$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F23").Value = 5
$ws.Range("I23").Value = 16
$ws.Range("K23").Value = 128.571428571429
$ws.Range("L23").Value = 33.333333333333
$ws.Range("M23").Value = 6.666666666666
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 12.037037037037
$ws.Range("I24").Value = 590
$ws.Range("J24").Value = 461
$ws.Range("K24").Value = 27.982646420824
$ws.Range("L24").Value = 44.607843137254
$ws.Range("M24").Value = 72.51461988304
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = 94.736842105263
$ws.Range("F25").Value = 103
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = 63.492063492063
$ws.Range("I25").Value = 493
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 64.333333333333
$ws.Range("L25").Value = 90.34749034749
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 75
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 22.222222222222
$ws.Range("I26").Value = 114
$ws.Range("J26").Value = 96
$ws.Range("K26").Value = 18.75
$ws.Range("L26").Value = 48.051948051948
$ws.Range("M26").Value = -2.564102564102
$ws.Range("D27").Value = 1
$ws.Range("I27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = -100
$ws.Range("K27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G27").Value = 1
$ws.Range("I27").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H27").Value = -100
$ws.Range("K27").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -63.636363636363
$ws.Range("L27").Value = -42.857142857142
$ws.Range("C28").Value = 3
$ws.Range("F28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D28").Value = "'0"
$ws.Range("A28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = "'***.*"
$ws.Range("A28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 11
$ws.Range("K28").Value = -21.428571428571
$ws.Range("L28").Value = 10
$ws.Range("N29").Value = -90.47619047619
$ws.Range("N30").Value = -89.473684210526
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 200
$ws.Range("L31").Value = -40
